$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 0.8234243614228944
$ws.Cells.Item(2, 3).Value2 = 0.2390411731991264
$ws.Cells.Item(2, 4).Value2 = 0.01114603848917994
$ws.Cells.Item(2, 6).Value2 = 0.4783276659663898
$ws.Cells.Item(2, 7).Value2 = 0.3296655200874881
$ws.Cells.Item(2, 8).Value2 = 0.4595714176470125
$ws.Cells.Item(2, 9).Value2 = 0.3002845041730531
$ws.Cells.Item(2, 14).Value2 = 0.8190571190906795

$ws.Cells.Item(3, 2).Value2 = 0.7214998676118967
$ws.Cells.Item(3, 3).Value2 = 0.2092395276970649
$ws.Cells.Item(3, 4).Value2 = 0.009965053555823999
$ws.Cells.Item(3, 6).Value2 = 0.4689669770731228
$ws.Cells.Item(3, 7).Value2 = 0.320977232158981
$ws.Cells.Item(3, 8).Value2 = 0.460325805726427
$ws.Cells.Item(3, 9).Value2 = 0.3038796859177353
$ws.Cells.Item(3, 14).Value2 = 0.8220607328432337

$ws.Cells.Item(4, 2).Value2 = 0.6588221578764717
$ws.Cells.Item(4, 3).Value2 = 0.1908649908417601
$ws.Cells.Item(4, 4).Value2 = 0.009236407118688561
$ws.Cells.Item(4, 6).Value2 = 0.4636780216299954
$ws.Cells.Item(4, 7).Value2 = 0.3160259442002484
$ws.Cells.Item(4, 8).Value2 = 0.4611342259807572
$ws.Cells.Item(4, 9).Value2 = 0.3063864754026007
$ws.Cells.Item(4, 14).Value2 = 0.8242948709632074

$ws.Cells.Item(5, 2).Value2 = 0.6332575624831804
$ws.Cells.Item(5, 3).Value2 = 0.1833582503887783
$ws.Cells.Item(5, 4).Value2 = 0.008938616136124011
$ws.Cells.Item(5, 6).Value2 = 0.4616375848263203
$ws.Cells.Item(5, 7).Value2 = 0.3141040627619418
$ws.Cells.Item(5, 8).Value2 = 0.4615502990809119
$ws.Cells.Item(5, 9).Value2 = 0.3074830104067452
$ws.Cells.Item(5, 14).Value2 = 0.8253034834620152

$ws.Cells.Item(6, 2).Value2 = 0.6290112301303736
$ws.Cells.Item(6, 3).Value2 = 0.1821106220582749
$ws.Cells.Item(6, 4).Value2 = 0.008889116799458918
$ws.Cells.Item(6, 6).Value2 = 0.4613056956842456
$ws.Cells.Item(6, 7).Value2 = 0.3137907051563076
$ws.Cells.Item(6, 8).Value2 = 0.4616246151814209
$ws.Cells.Item(6, 9).Value2 = 0.3076696109341626
$ws.Cells.Item(6, 14).Value2 = 0.8254768969764399

$ws.Cells.Item(7, 2).Value2 = 0.6584774760963228
$ws.Cells.Item(7, 3).Value2 = 0.1907638288517433
$ws.Cells.Item(7, 4).Value2 = 0.009232394464234517
$ws.Cells.Item(7, 6).Value2 = 0.4636500391369012
$ws.Cells.Item(7, 7).Value2 = 0.3159996378911529
$ws.Cells.Item(7, 8).Value2 = 0.4611394867367977
$ws.Cells.Item(7, 9).Value2 = 0.3064009603547255
$ws.Cells.Item(7, 14).Value2 = 0.8243080757400563

$ws.Cells.Item(8, 2).Value2 = 0.7883012982027253
$ws.Cells.Item(8, 3).Value2 = 0.2287814904606478
$ws.Cells.Item(8, 4).Value2 = 0.01073957805147785
$ws.Cells.Item(8, 6).Value2 = 0.4750046546502773
$ws.Cells.Item(8, 7).Value2 = 0.3265898749196054
$ws.Cells.Item(8, 8).Value2 = 0.4597597736787407
$ws.Cells.Item(8, 9).Value2 = 0.3014618326369245
$ws.Cells.Item(8, 14).Value2 = 0.8200119195486195

$ws.Cells.Item(9, 2).Value2 = 1.04208900937823
$ws.Cells.Item(9, 3).Value2 = 0.3027243139620452
$ws.Cells.Item(9, 4).Value2 = 0.01366639802007086
$ws.Cells.Item(9, 6).Value2 = 0.5009315587484622
$ws.Cells.Item(9, 7).Value2 = 0.3504271575938844
$ws.Cells.Item(9, 8).Value2 = 0.4598018282496952
$ws.Cells.Item(9, 9).Value2 = 0.294163625853443
$ws.Cells.Item(9, 14).Value2 = 0.814675716803805

$ws.Cells.Item(10, 2).Value2 = 1.228029407099825
$ws.Cells.Item(10, 3).Value2 = 0.3566768823969824
$ws.Cells.Item(10, 4).Value2 = 0.01579825279128499
$ws.Cells.Item(10, 6).Value2 = 0.522245074237091
$ws.Cells.Item(10, 7).Value2 = 0.3698529550440668
$ws.Cells.Item(10, 8).Value2 = 0.4615204330958989
$ws.Cells.Item(10, 9).Value2 = 0.2902742258208022
$ws.Cells.Item(10, 14).Value2 = 0.8126323559473292

$ws.Cells.Item(11, 2).Value2 = 1.312500563948618
$ws.Cells.Item(11, 3).Value2 = 0.3811402683526808
$ws.Cells.Item(11, 4).Value2 = 0.01676389016963498
$ws.Cells.Item(11, 6).Value2 = 0.5324403204424186
$ws.Cells.Item(11, 7).Value2 = 0.3791144895305223
$ws.Cells.Item(11, 8).Value2 = 0.4626715298432487
$ws.Cells.Item(11, 9).Value2 = 0.2888282964741897
$ws.Cells.Item(11, 14).Value2 = 0.8121093590750235

$ws.Cells.Item(12, 2).Value2 = 1.344470306865276
$ws.Cells.Item(12, 3).Value2 = 0.3903923048377465
$ws.Cells.Item(12, 4).Value2 = 0.01712893466666543
$ws.Cells.Item(12, 6).Value2 = 0.5363733735614886
$ws.Cells.Item(12, 7).Value2 = 0.3826833276893353
$ws.Cells.Item(12, 8).Value2 = 0.4631607347953519
$ws.Cells.Item(12, 9).Value2 = 0.288327562248373
$ws.Cells.Item(12, 14).Value2 = 0.811969683973885

$ws.Cells.Item(13, 2).Value2 = 1.337585855576151
$ws.Cells.Item(13, 3).Value2 = 0.3884002379398339
$ws.Cells.Item(13, 4).Value2 = 0.01705034382197113
$ws.Cells.Item(13, 6).Value2 = 0.5355230953093866
$ws.Cells.Item(13, 7).Value2 = 0.3819119605190764
$ws.Cells.Item(13, 8).Value2 = 0.4630530015227805
$ws.Cells.Item(13, 9).Value2 = 0.2884333176540146
$ws.Cells.Item(13, 14).Value2 = 0.8119971707014457

$ws.Cells.Item(14, 2).Value2 = 1.315131094826086
$ws.Cells.Item(14, 3).Value2 = 0.3819016750770743
$ws.Cells.Item(14, 4).Value2 = 0.01679393519172834
$ws.Cells.Item(14, 6).Value2 = 0.5327624425936932
$ws.Cells.Item(14, 7).Value2 = 0.3794068598903522
$ws.Cells.Item(14, 8).Value2 = 0.4627107074479397
$ws.Cells.Item(14, 9).Value2 = 0.2887861609146896
$ws.Cells.Item(14, 14).Value2 = 0.8120966985284355

$ws.Cells.Item(15, 2).Value2 = 1.301374569848804
$ws.Cells.Item(15, 3).Value2 = 0.377919585106838
$ws.Cells.Item(15, 4).Value2 = 0.01663679586290101
$ws.Cells.Item(15, 6).Value2 = 0.5310808977181125
$ws.Cells.Item(15, 7).Value2 = 0.3778804668392297
$ws.Cells.Item(15, 8).Value2 = 0.4625079910276497
$ws.Cells.Item(15, 9).Value2 = 0.289008392453983
$ws.Cells.Item(15, 14).Value2 = 0.8121652615242141

$ws.Cells.Item(16, 2).Value2 = 1.222506616238945
$ws.Cells.Item(16, 3).Value2 = 0.3550765152811266
$ws.Cells.Item(16, 4).Value2 = 0.01573506046502615
$ws.Cells.Item(16, 6).Value2 = 0.521588886460691
$ws.Cells.Item(16, 7).Value2 = 0.3692562910051294
$ws.Cells.Item(16, 8).Value2 = 0.4614526549286069
$ws.Cells.Item(16, 9).Value2 = 0.2903752541000379
$ws.Cells.Item(16, 14).Value2 = 0.8126747080203103

$ws.Cells.Item(17, 2).Value2 = 1.174093576873133
$ws.Cells.Item(17, 3).Value2 = 0.3410424042062346
$ws.Cells.Item(17, 4).Value2 = 0.01518079387032856
$ws.Cells.Item(17, 6).Value2 = 0.515894158716506
$ws.Cells.Item(17, 7).Value2 = 0.3640748502167099
$ws.Cells.Item(17, 8).Value2 = 0.4608999717792557
$ws.Cells.Item(17, 9).Value2 = 0.2912968224636501
$ws.Cells.Item(17, 14).Value2 = 0.8130912972097519

$ws.Cells.Item(18, 2).Value2 = 1.14623701531724
$ws.Cells.Item(18, 3).Value2 = 0.3329628468440831
$ws.Cells.Item(18, 4).Value2 = 0.01486160481321974
$ws.Cells.Item(18, 6).Value2 = 0.512665699422783
$ws.Cells.Item(18, 7).Value2 = 0.3611345540391682
$ws.Cells.Item(18, 8).Value2 = 0.4606168305178926
$ws.Cells.Item(18, 9).Value2 = 0.2918573037892394
$ws.Cells.Item(18, 14).Value2 = 0.8133691784282462

$ws.Cells.Item(19, 2).Value2 = 1.136803469180052
$ws.Cells.Item(19, 3).Value2 = 0.3302259649711345
$ws.Cells.Item(19, 4).Value2 = 0.01475346671015387
$ws.Cells.Item(19, 6).Value2 = 0.5115806555054689
$ws.Cells.Item(19, 7).Value2 = 0.360145859244497
$ws.Cells.Item(19, 8).Value2 = 0.4605269254096527
$ws.Cells.Item(19, 9).Value2 = 0.2920522883892431
$ws.Cells.Item(19, 14).Value2 = 0.8134698399241103

$ws.Cells.Item(20, 2).Value2 = 1.179248340078004
$ws.Cells.Item(20, 3).Value2 = 0.3425371379676676
$ws.Cells.Item(20, 4).Value2 = 0.01523983698849207
$ws.Cells.Item(20, 6).Value2 = 0.5164955045920152
$ws.Cells.Item(20, 7).Value2 = 0.3646222860823514
$ws.Cells.Item(20, 8).Value2 = 0.4609552081215185
$ws.Cells.Item(20, 9).Value2 = 0.2911955694045645
$ws.Cells.Item(20, 14).Value2 = 0.8130429904120575

$ws.Cells.Item(21, 2).Value2 = 1.321727091338687
$ws.Cells.Item(21, 3).Value2 = 0.3838107803078969
$ws.Cells.Item(21, 4).Value2 = 0.01686926567759173
$ws.Cells.Item(21, 6).Value2 = 0.5335713461754921
$ws.Cells.Item(21, 7).Value2 = 0.3801409892163576
$ws.Cells.Item(21, 8).Value2 = 0.462809799068836
$ws.Cells.Item(21, 9).Value2 = 0.2886812495032203
$ws.Cells.Item(21, 14).Value2 = 0.8120658812474915

$ws.Cells.Item(22, 2).Value2 = 1.414741732651294
$ws.Cells.Item(22, 3).Value2 = 0.4107170925403807
$ws.Cells.Item(22, 4).Value2 = 0.01793055834883717
$ws.Cells.Item(22, 6).Value2 = 0.5451532440643803
$ws.Cells.Item(22, 7).Value2 = 0.3906432636472488
$ws.Cells.Item(22, 8).Value2 = 0.4643327061128133
$ws.Cells.Item(22, 9).Value2 = 0.2873109700939622
$ws.Cells.Item(22, 14).Value2 = 0.8117674869457403

$ws.Cells.Item(23, 2).Value2 = 1.365107945887814
$ws.Cells.Item(23, 3).Value2 = 0.3963630202895843
$ws.Cells.Item(23, 4).Value2 = 0.01736446682171788
$ws.Cells.Item(23, 6).Value2 = 0.538933006470657
$ws.Cells.Item(23, 7).Value2 = 0.3850048619195405
$ws.Cells.Item(23, 8).Value2 = 0.4634913935387743
$ws.Cells.Item(23, 9).Value2 = 0.2880172356561914
$ws.Cells.Item(23, 14).Value2 = 0.8118956443212966

$ws.Cells.Item(24, 2).Value2 = 1.176917943994113
$ws.Cells.Item(24, 3).Value2 = 0.341861403583664
$ws.Cells.Item(24, 4).Value2 = 0.01521314525546558
$ws.Cells.Item(24, 6).Value2 = 0.5162234943808528
$ws.Cells.Item(24, 7).Value2 = 0.364374670192646
$ws.Cells.Item(24, 8).Value2 = 0.4609301279958231
$ws.Cells.Item(24, 9).Value2 = 0.2912412504321935
$ws.Cells.Item(24, 14).Value2 = 0.8130647103585318

$ws.Cells.Item(25, 2).Value2 = 0.9735212901853743
$ws.Cells.Item(25, 3).Value2 = 0.2827862561440497
$ws.Cells.Item(25, 4).Value2 = 0.01287779329090455
$ws.Cells.Item(25, 6).Value2 = 0.4935222974732341
$ws.Cells.Item(25, 7).Value2 = 0.3436456421788563
$ws.Cells.Item(25, 8).Value2 = 0.4594949878401451
$ws.Cells.Item(25, 9).Value2 = 0.2958806213648835
$ws.Cells.Item(25, 14).Value2 = 0.8157892647763134
